$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (number format / style / boolean type) of the previous
# trade row (row 3) down onto the new row 4 so the new row inherits the same
# date style on column A/G and boolean cell types on B/G/I.
$ws.Range("A3:I3").Copy($ws.Range("A4:I4"))

# Now overwrite the copied values with this new trade's data (20 minute trade).
$ws.Range("A4").Value = 42641.545856481483
$ws.Range("B4").Value = $true
$ws.Range("C4").Value = 10001.969999999999
$ws.Range("D4").Value = 9985
$ws.Range("E4").Value = 309.77999999999997
$ws.Range("F4").Value = 308.73
$ws.Range("G4").Value = $true
$ws.Range("H4").Value = -0.34
$ws.Range("I4").Value = $false

# Column C ("Principle") needs to widen slightly to fit the new, wider value
# (10001.97 vs the previous 9985/10000), matching the workbook's bestFit
# auto-sizing behaviour.
$ws.Columns.Item(3).ColumnWidth = 8.1
